$wb = $excel.ActiveWorkbook

# --- 1. BANK0x3D: remove "viewTab" and "loadedView" rows (now calculated automatically) ---
$bank3d = $wb.Worksheets.Item("BANK0x3D")
$bank3d.Range("A4:E5").ClearContents()
$bank3d.Range("A3").Select()

# --- 2. Add new "Sprite Addresses" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$spriteSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$spriteSheet.Name = "Sprite Addresses"

$spriteSheet.Cells.Item(1, 1).Value2 = "Index"
$spriteSheet.Cells.Item(1, 2).Value2 = "Value"

$baseAddr = 0xEA00
for ($i = 0; $i -lt 135; $i++) {
    $addr = $baseAddr + $i * 0x200
    $hexText = "0X" + [System.Convert]::ToString($addr, 16).ToUpper()
    $spriteSheet.Cells.Item($i + 2, 1).Value2 = $i
    $spriteSheet.Cells.Item($i + 2, 2).Value2 = $hexText
}

$spriteSheet.Cells.Item(1, 4).Value2 = "All possible sprite addresses and their index in _bESpriteAllocTable"

$spriteSheet.Columns.Item(4).ColumnWidth = 119.18

$spriteSheet.Activate()
$spriteSheet.Range("D1").Select()
$excel.ActiveWindow.ScrollRow = 4

# --- 3. Restore the active tab to BANK0x3C ---
$bank3c = $wb.Worksheets.Item("BANK0x3C")
$bank3c.Activate()
